$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.38"
$ws.Range("E2").Value = "'2.61%"
$ws.Range("D3").Value = "'35.96"
$ws.Range("E3").Value = "'1.60%"
$ws.Range("D4").Value = "'5.066"
$ws.Range("E4").Value = "'0.52%"
$ws.Range("D5").Value = "'0.08136"
$ws.Range("E5").Value = "'2.19%"
$ws.Range("D6").Value = "'1.947"
$ws.Range("E6").Value = "'2.36%"
$ws.Range("D7").Value = "'4.138"
$ws.Range("E7").Value = "'1.99%"
$ws.Range("D8").Value = "'7.800"
$ws.Range("E8").Value = "'0.10%"
$ws.Range("D9").Value = "'0.9361"
$ws.Range("E9").Value = "'1.52%"
$ws.Range("D10").Value = "'0.1324"
$ws.Range("E10").Value = "'-6.53%"
$ws.Range("D11").Value = "'0.1919"
$ws.Range("E11").Value = "'1.25%"
$ws.Range("D12").Value = "'0.09266"
$ws.Range("E12").Value = "'1.13%"
$ws.Range("D13").Value = "'0.03509"
$ws.Range("E13").Value = "'2.79%"
$ws.Range("D14").Value = "'0.09883"
$ws.Range("E14").Value = "'0.24%"
$ws.Range("D15").Value = "'0.001444"
$ws.Range("E15").Value = "'3.51%"
$ws.Range("D16").Value = "'0.005749"
$ws.Range("E16").Value = "'-1.37%"
$ws.Range("D17").Value = "'3.605"
$ws.Range("E17").Value = "'2.51%"
$ws.Range("E18").Value = "'0.01%"
$ws.Range("E19").Value = "'0.77%"
$ws.Range("E20").Value = "'3.66%"
$ws.Range("D21").Value = "'5.174"
$ws.Range("E21").Value = "'2.02%"
$ws.Range("D22").Value = "'0.2614"
$ws.Range("E22").Value = "'8.62%"
$ws.Range("D23").Value = "'0.04386"
$ws.Range("E23").Value = "'-2.33%"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'0.27%"
$ws.Range("D25").Value = "'0.004779"
$ws.Range("E25").Value = "'0.26%"
$ws.Range("E26").Value = "'33.76%"
$ws.Range("D27").Value = "'0.0003125"
$ws.Range("E27").Value = "'3.84%"
$ws.Range("D39").Value = "'0.01997"
$ws.Range("E39").Value = "'5.44%"
$ws.Range("D40").Value = "'0.05059"
$ws.Range("E40").Value = "'7.30%"
$ws.Range("D41").Value = "'0.01123"
$ws.Range("E41").Value = "'15.43%"
$ws.Range("D42").Value = "'0.007622"
$ws.Range("E42").Value = "'3.18%"
$ws.Range("D43").Value = "'0.1382"
$ws.Range("E43").Value = "'4.40%"
$ws.Range("D44").Value = "'0.002097"
$ws.Range("D45").Value = "'0.01129"
$ws.Range("E45").Value = "'8.79%"
$ws.Range("D46").Value = "'0.00006393"
$ws.Range("E46").Value = "'1.92%"
$ws.Range("E47").Value = "'-0.43%"
$ws.Range("D48").Value = "'64.96"
$ws.Range("E48").Value = "'0.75%"
$ws.Range("D49").Value = "'0.001189"
$ws.Range("E49").Value = "'-28.54%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.43%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.43%"

Write-Host "Applied all cell updates"
